$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Extend the "Rockwell" specs text (E3) with additional spec lines ---
$ws.Range("E3").Value = "processeur : de 350 MHz à 1 GHz`nRAM: 256 MB `nDiagonales : 5,7 pouces`nrésolution : 320 x 240px`nDimensions (L x H x P) : 185 x 152 x 68mm`nCouleurs : N/A`nType d'affichage : TFT`nPoids : 594`nProtection :`tIP66 (face avant)"

# --- 2) Fill in the previously empty "Siemens" specs cell (E2) with a new text ---
$ws.Range("E2").Value = "processeur : MRA`nFlash / RAM : 512 kbyte`nDiagonales : 5,7 pouces`nrésolution : 320 x 240px`nDimensions (L x H x P) : 212 x 156 x 44 mm `nCouleurs : 2`nType d'affichage : LCD`nPoids : 0,75 kg`nProtection : IP65 (face avant) / IP20 (panneau arriere)"

# --- 3) Harmonise the "Caractéristiques" column formatting (wrap text like E4:E6) ---
$ws.Range("E2:E6").WrapText = $true

# --- 4) Harmonise the "Prix (HT)" column number format (plain integer, no currency) ---
$ws.Range("F2:F6").NumberFormat = "0"

# --- 5) Widen column F slightly ---
$ws.Columns("F").ColumnWidth = 15.5546875

# --- 6) Update the view: scroll so row 2 is at top, and select E2 ---
$ws.Range("E2").Select()
$excel.ActiveWindow.ScrollRow = 2
